# Refresh the cryptocurrency ranking table (rows 2-51) with the latest
# price/volume snapshot. Column A (rank) and the header row are untouched.
#
# Numeric-looking price strings (e.g. "1.00", "40.40") are written with a
# leading apostrophe so Excel keeps them as literal text (matching the
# source sheet, which stores every Price/Volume cell as text) instead of
# normalizing them into plain numbers (which would turn "1.00" into "1").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "48.359.28"
$ws.Range("E2").Value = "  +1.98%  "
# Row 3: Ethereum
$ws.Range("D3").Value = "2.519.49"
$ws.Range("E3").Value = "  +0.35%  "
# Row 4: TetherUSD
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "
# Row 5: BNB
$ws.Range("D5").Value = "'323.57"
$ws.Range("E5").Value = "  -0.21%  "
# Row 6: Solana
$ws.Range("D6").Value = "'109.39"
$ws.Range("E6").Value = "  -0.31%  "
# Row 7: XRP
$ws.Range("E7").Value = "  -0.24%  "
# Row 8: USDC
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.03%  "
# Row 9: Cardano
$ws.Range("D9").Value = "'0.562"
$ws.Range("E9").Value = "  +3.83%  "
# Row 10: Avalanche
$ws.Range("D10").Value = "'40.40"
$ws.Range("E10").Value = "  +2.89%  "
# Row 11: Dogecoin
$ws.Range("D11").Value = "'0.0822"
$ws.Range("E11").Value = "  +0.31%  "
# Row 12: Chainlink
$ws.Range("D12").Value = "'19.62"
$ws.Range("E12").Value = "  +5.34%  "
# Row 13: TRON
$ws.Range("E13").Value = "  +0.74%  "
# Row 14: Polkadot
$ws.Range("E14").Value = "  -0.33%  "
# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.907.28"
$ws.Range("E15").Value = "  +0.08%  "
# Row 16: WrappedEther
$ws.Range("D16").Value = "2.517.24"
$ws.Range("E16").Value = "  +0.27%  "
# Row 17: Polygon
$ws.Range("D17").Value = "'0.856"
$ws.Range("E17").Value = "  -0.78%  "
# Row 18: WrappedBTC
$ws.Range("D18").Value = "48.216.20"
$ws.Range("E18").Value = "  +1.73%  "
# Row 19: InternetComputer(DFINITY)
$ws.Range("D19").Value = "'13.46"
$ws.Range("E19").Value = "  +3.79%  "
# Row 20: Uniswap
$ws.Range("D20").Value = "'6.65"
$ws.Range("E20").Value = "  -1.03%  "
# Row 21: ShibaInu
$ws.Range("D21").Value = "0.0₃0947"
$ws.Range("E21").Value = "  -0.15%  "
# Row 22: ImmutableX
$ws.Range("D22").Value = "'2.74"
$ws.Range("E22").Value = "  +2.56%  "
# Row 23: Litecoin
$ws.Range("D23").Value = "'72.30"
$ws.Range("E23").Value = "  +1.94%  "
# Row 24: BitcoinCash
$ws.Range("D24").Value = "'268.16"
$ws.Range("E24").Value = "  +7.19%  "
# Row 25: PancakeSwap
$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = "  -1.59%  "
# Row 26: EthereumClassic
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'26.21"
$ws.Range("E26").Value = "  -0.13%  "
# Row 27: Dai
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.04%  "
# Row 28: Cosmos
$ws.Range("E28").Value = "  +1.56%  "
# Row 29: Kaspa
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.143"
$ws.Range("E29").Value = "  +5.18%  "
# Row 30: Toncoin
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.20"
$ws.Range("E30").Value = "  -4.23%  "
# Row 31: InjectiveProtocol
$ws.Range("D31").Value = "'35.38"
$ws.Range("E31").Value = "  -1.11%  "
# Row 32: OKB
$ws.Range("D32").Value = "'49.96"
$ws.Range("E32").Value = "  -0.76%  "
# Row 33: Celestia
$ws.Range("D33").Value = "'20.06"
$ws.Range("E33").Value = "  +0.53%  "
# Row 34: Filecoin
$ws.Range("E34").Value = "  -0.91%  "
# Row 36: Hedera
$ws.Range("D36").Value = "'0.0788"
$ws.Range("E36").Value = "  -1.16%  "
# Row 37: ARBITRUM
$ws.Range("E37").Value = "  -1.34%  "
# Row 38: RenderToken
$ws.Range("D38").Value = "'4.72"
$ws.Range("E38").Value = "  -0.48%  "
# Row 39: LidoDAOToken
$ws.Range("E39").Value = "  -1.05%  "
# Row 40: Stellar
$ws.Range("E40").Value = "  -0.27%  "
# Row 41: EnergySwap
$ws.Range("D41").Value = "'22.45"
$ws.Range("E41").Value = "  +4.95%  "
# Row 42: Monero
$ws.Range("D42").Value = "'119.05"
$ws.Range("E42").Value = "  -3.09%  "
# Row 43: WEMIXToken
$ws.Range("D43").Value = "'2.18"
$ws.Range("E43").Value = "  -3.44%  "
# Row 44: VeChain
$ws.Range("E44").Value = "  +0.34%  "
# Row 45: Maker
$ws.Range("D45").Value = "2.000.97"
$ws.Range("E45").Value = "  -0.07%  "
# Row 46: NEARProtocol
$ws.Range("E46").Value = "  +0.64%  "
# Row 47: ApeXProtocol
$ws.Range("E47").Value = "  -2.96%  "
# Row 48: Stacks
$ws.Range("E48").Value = "  +3.65%  "
# Row 49: FraxShare
$ws.Range("D49").Value = "'9.10"
$ws.Range("E49").Value = "  +0.23%  "
# Row 50: THORChain
$ws.Range("D50").Value = "'5.25"
$ws.Range("E50").Value = "  -1.07%  "
# Row 51: BitcoinSV
$ws.Range("D51").Value = "'80.15"
$ws.Range("E51").Value = "  +1.95%  "
